$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H47").Value = 6500
$ws.Range("I47").Value = 6500
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 6500
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -5528

$ws = $wb.Worksheets("ALC")
$ws.Range("H98").Value = 4918.577
$ws.Range("I98").Value = 4918.577
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 4918.577
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -3420.577
$ws.Range("N98").ClearContents()

$ws = $wb.Worksheets("ALC")
$ws.Range("H106").Value = 14737158
$ws.Range("I106").Value = 33979.582
$ws.Range("J106").Value = 166670000
$ws.Range("K106").Value = 33979.582
$ws.Range("L106").Value = 166670000
$ws.Range("M106").Value = -33348.582
$ws.Range("N106").Value = -166671262

$ws = $wb.Worksheets("ALC")
$ws.Range("H116").Value = 5850
$ws.Range("I116").Value = 1500
$ws.Range("J116").Value = 7300
$ws.Range("K116").Value = 1500
$ws.Range("L116").Value = 7300
$ws.Range("M116").Value = 1942
$ws.Range("N116").Value = -14184

$ws = $wb.Worksheets("ALC")
$ws.Range("H122").Value = 4918.577
$ws.Range("I122").Value = 4918.577
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 14755.731
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -12305.731
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets("ALC")
$ws.Range("H135").Value = 5340.207
$ws.Range("I135").Value = 719.2381
$ws.Range("J135").Value = 17470.25
$ws.Range("K135").Value = 6473.142900000001
$ws.Range("L135").Value = 157232.25
$ws.Range("M135").Value = -3938.142900000001
$ws.Range("N135").Value = -162302.25

$ws = $wb.Worksheets("ARM")
$ws.Range("H2").Value = 41667372
$ws.Range("I2").Value = 66667170
$ws.Range("J2").Value = 1044.8889
$ws.Range("K2").Value = 66667170
$ws.Range("L2").Value = 1044.8889
$ws.Range("M2").Value = -66667057
$ws.Range("N2").Value = -1270.8889

$ws = $wb.Worksheets("ARM")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()

$ws = $wb.Worksheets("ARM")
$ws.Range("H32").Value = 18092.242
$ws.Range("I32").Value = 18144.658
$ws.Range("J32").Value = 17882.578
$ws.Range("K32").Value = 18144.658
$ws.Range("L32").Value = 17882.578
$ws.Range("M32").Value = -17857.658
$ws.Range("N32").Value = -18456.578

$ws = $wb.Worksheets("ARM")
$ws.Range("H45").Value = 20834262
$ws.Range("I45").Value = 33334196
$ws.Range("J45").Value = 1038
$ws.Range("K45").Value = 33334196
$ws.Range("L45").Value = 1038
$ws.Range("M45").Value = -33333819
$ws.Range("N45").Value = -1792

$ws = $wb.Worksheets("ARM")
$ws.Range("H59").Value = 10000
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 10000
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 10000
$ws.Range("N59").Value = -11608

$ws = $wb.Worksheets("ARM")
$ws.Range("H116").Value = 41667372
$ws.Range("I116").Value = 66667170
$ws.Range("J116").Value = 1044.8889
$ws.Range("K116").Value = 66667170
$ws.Range("L116").Value = 1044.8889
$ws.Range("M116").Value = -66664876
$ws.Range("N116").Value = -5632.8889

$ws = $wb.Worksheets("BSM")
$ws.Range("H3").Value = 41667372
$ws.Range("I3").Value = 66667170
$ws.Range("J3").Value = 1044.8889
$ws.Range("K3").Value = 66667170
$ws.Range("L3").Value = 1044.8889
$ws.Range("M3").Value = -66667056
$ws.Range("N3").Value = -1272.8889

$ws = $wb.Worksheets("BSM")
$ws.Range("H86").Value = 2100
$ws.Range("I86").Value = 2107.6924
$ws.Range("J86").Value = 2066.6667
$ws.Range("K86").Value = 2107.6924
$ws.Range("L86").Value = 2066.6667
$ws.Range("M86").Value = -984.6923999999999
$ws.Range("N86").Value = -4312.6667

$ws = $wb.Worksheets("BSM")
$ws.Range("H89").Value = 2100
$ws.Range("I89").Value = 2107.6924
$ws.Range("J89").Value = 2066.6667
$ws.Range("K89").Value = 10538.462
$ws.Range("L89").Value = 10333.3335
$ws.Range("M89").Value = -4922.462
$ws.Range("N89").Value = -21565.3335

$ws = $wb.Worksheets("CRP")
$ws.Range("H31").Value = 2655.7568
$ws.Range("I31").Value = 2353.1
$ws.Range("J31").Value = 3952.8572
$ws.Range("K31").Value = 2353.1
$ws.Range("L31").Value = 3952.8572
$ws.Range("M31").Value = -2058.1
$ws.Range("N31").Value = -4542.8572

$ws = $wb.Worksheets("CRP")
$ws.Range("H34").Value = 2655.7568
$ws.Range("I34").Value = 2353.1
$ws.Range("J34").Value = 3952.8572
$ws.Range("K34").Value = 2353.1
$ws.Range("L34").Value = 3952.8572
$ws.Range("M34").Value = -2151.1
$ws.Range("N34").Value = -4356.8572

$ws = $wb.Worksheets("CUL")
$ws.Range("H11").Value = 9093345
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 9093345
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 27280035
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -27280315

$ws = $wb.Worksheets("CUL")
$ws.Range("H26").Value = 10527289
$ws.Range("I26").Value = 60
$ws.Range("J26").Value = 13334550
$ws.Range("K26").Value = 180
$ws.Range("L26").Value = 40003650
$ws.Range("M26").Value = 108
$ws.Range("N26").Value = -40004226

$ws = $wb.Worksheets("CUL")
$ws.Range("H34").Value = 1245.5385
$ws.Range("I34").Value = 372.25
$ws.Range("J34").Value = 1633.6666
$ws.Range("K34").Value = 1116.75
$ws.Range("L34").Value = 4900.9998
$ws.Range("M34").Value = -1032.75
$ws.Range("N34").Value = -5068.9998

$ws = $wb.Worksheets("CUL")
$ws.Range("H39").Value = 3383.4614
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 3383.4614
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 10150.3842
$ws.Range("N39").Value = -10738.3842

$ws = $wb.Worksheets("CUL")
$ws.Range("H51").Value = 800
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 800
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 2400
$ws.Range("N51").Value = -3320

$ws = $wb.Worksheets("CUL")
$ws.Range("H52").Value = 544.3333
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 544.3333
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 1632.9999
$ws.Range("N52").Value = -2164.9999

$ws = $wb.Worksheets("CUL")
$ws.Range("H62").Value = 3500
$ws.Range("I62").Value = 2500
$ws.Range("J62").Value = 4500
$ws.Range("K62").Value = 7500
$ws.Range("L62").Value = 13500
$ws.Range("M62").Value = -6814
$ws.Range("N62").Value = -14872

$ws = $wb.Worksheets("CUL")
$ws.Range("H65").Value = 3500
$ws.Range("I65").Value = 2500
$ws.Range("J65").Value = 4500
$ws.Range("K65").Value = 22500
$ws.Range("L65").Value = 40500
$ws.Range("M65").Value = -19068
$ws.Range("N65").Value = -47364

$ws = $wb.Worksheets("CUL")
$ws.Range("H80").Value = 2812.875
$ws.Range("I80").Value = 950.75
$ws.Range("J80").Value = 4675
$ws.Range("K80").Value = 2852.25
$ws.Range("L80").Value = 14025
$ws.Range("M80").Value = -1916.25
$ws.Range("N80").Value = -15897

$ws = $wb.Worksheets("CUL")
$ws.Range("H82").Value = 3980
$ws.Range("I82").Value = 3300
$ws.Range("J82").Value = 5000
$ws.Range("K82").Value = 9900
$ws.Range("L82").Value = 15000
$ws.Range("M82").Value = -9494
$ws.Range("N82").Value = -15812

$ws = $wb.Worksheets("CUL")
$ws.Range("H83").Value = 2812.875
$ws.Range("I83").Value = 950.75
$ws.Range("J83").Value = 4675
$ws.Range("K83").Value = 8556.75
$ws.Range("L83").Value = 42075
$ws.Range("M83").Value = -3876.75
$ws.Range("N83").Value = -51435

$ws = $wb.Worksheets("CUL")
$ws.Range("H85").Value = 3980
$ws.Range("I85").Value = 3300
$ws.Range("J85").Value = 5000
$ws.Range("K85").Value = 9900
$ws.Range("L85").Value = 15000
$ws.Range("M85").Value = -8496
$ws.Range("N85").Value = -17808

$ws = $wb.Worksheets("CUL")
$ws.Range("H116").Value = 9845.23
$ws.Range("I116").Value = 12998.777
$ws.Range("J116").Value = 2749.75
$ws.Range("K116").Value = 38996.331
$ws.Range("L116").Value = 8249.25
$ws.Range("M116").Value = -35554.331
$ws.Range("N116").Value = -15133.25

$ws = $wb.Worksheets("CUL")
$ws.Range("H130").Value = 2472.5
$ws.Range("I130").Value = 2390
$ws.Range("J130").Value = 2500
$ws.Range("K130").Value = 7170
$ws.Range("L130").Value = 7500
$ws.Range("M130").Value = -2150
$ws.Range("N130").Value = -17540

$ws = $wb.Worksheets("CUL")
$ws.Range("H131").Value = 738.55
$ws.Range("I131").Value = 327.33334
$ws.Range("J131").Value = 794.625
$ws.Range("K131").Value = 982.0000200000001
$ws.Range("L131").Value = 2383.875
$ws.Range("M131").Value = 4057.99998
$ws.Range("N131").Value = -12463.875

$ws = $wb.Worksheets("GSM")
$ws.Range("H113").Value = 27779010
$ws.Range("I113").Value = 125000600
$ws.Range("J113").Value = 1412.7142
$ws.Range("K113").Value = 125000600
$ws.Range("L113").Value = 1412.7142
$ws.Range("M113").Value = -124998430
$ws.Range("N113").Value = -5752.7142

$ws = $wb.Worksheets("GSM")
$ws.Range("H126").Value = 2087.5
$ws.Range("I126").Value = 1663.1578
$ws.Range("J126").Value = 3700
$ws.Range("K126").Value = 4989.4734
$ws.Range("L126").Value = 11100
$ws.Range("M126").Value = -2519.4734
$ws.Range("N126").Value = -16040

$ws = $wb.Worksheets("GSM")
$ws.Range("H132").Value = 11834.167
$ws.Range("I132").Value = 17644.572
$ws.Range("J132").Value = 3699.6
$ws.Range("K132").Value = 52933.716
$ws.Range("L132").Value = 11098.8
$ws.Range("M132").Value = -50403.716
$ws.Range("N132").Value = -16158.8

$ws = $wb.Worksheets("LTW")
$ws.Range("H14").Value = 1200
$ws.Range("I14").Value = 1200
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1200
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -1028

$ws = $wb.Worksheets("LTW")
$ws.Range("H82").Value = 1056
$ws.Range("I82").Value = 1004.2727
$ws.Range("J82").Value = 1625
$ws.Range("K82").Value = 1004.2727
$ws.Range("L82").Value = 1625
$ws.Range("M82").Value = -643.2727
$ws.Range("N82").Value = -2347

$ws = $wb.Worksheets("LTW")
$ws.Range("H85").Value = 1056
$ws.Range("I85").Value = 1004.2727
$ws.Range("J85").Value = 1625
$ws.Range("K85").Value = 1004.2727
$ws.Range("L85").Value = 1625
$ws.Range("M85").Value = 243.7273
$ws.Range("N85").Value = -4121

$ws = $wb.Worksheets("LTW")
$ws.Range("H132").Value = 9802.370000000001
$ws.Range("I132").Value = 14391.375
$ws.Range("J132").Value = 3127.4546
$ws.Range("K132").Value = 43174.125
$ws.Range("L132").Value = 9382.363799999999
$ws.Range("M132").Value = -40644.125
$ws.Range("N132").Value = -14442.3638

$ws = $wb.Worksheets("LTW")
$ws.Range("H136").Value = 3434.0728
$ws.Range("I136").Value = 3424.5625
$ws.Range("J136").Value = 3499.2856
$ws.Range("K136").Value = 10273.6875
$ws.Range("L136").Value = 10497.8568
$ws.Range("M136").Value = -7723.6875
$ws.Range("N136").Value = -15597.8568

$ws = $wb.Worksheets("WVR")
$ws.Range("H136").Value = 1760.3922
$ws.Range("I136").Value = 1672.7273
$ws.Range("J136").Value = 1921.1111
$ws.Range("K136").Value = 5018.1819
$ws.Range("L136").Value = 5763.3333
$ws.Range("M136").Value = -2468.1819
$ws.Range("N136").Value = -10863.3333

